$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing data rows 2..116 down to 3..117)
$ws.Rows.Item(2).Insert()

# Copy the number/cell formatting from row 3 (a normal data row) onto the
# newly inserted row 2 so it matches the rest of the table (date style on
# column D, plain/general style elsewhere) instead of inheriting the header
# row's bold style.
$ws.Range("A3:R3").Copy()
$ws.Range("A2:R2").PasteSpecial(-4122)

# Populate the new row with the new weekly price observation.
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44921
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 100114007
$ws.Range("G2").Value = "Jengibre"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 610
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("N2").Value = "`$/caja 13 kilos"
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 1154
$ws.Range("Q2").Value = 13
$ws.Range("R2").Value = "Hortaliza"
